$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Voice Lines - main")

# --- Trim / rewrite the "Comments" column (F) text on a handful of rows. ---
$ws.Range("F4").Value = "OPTION ""Go right""`n"
$ws.Range("F10").Value = "Comment for a line.`nAnother comment for the same line."
$ws.Range("F14").Value = "This is for all the barks.`n(1/7) "
$ws.Range("F16").Value = "(3/7) "

# --- Column F (Comments) got a bit narrower after the trims above. ---
$ws.Columns("F").ColumnWidth = 32.960625

# --- Row heights for rows 4 and 10 shrank to fit the now-shorter wrapped text. ---
$ws.Rows("4").RowHeight = 13.5
$ws.Rows("10").RowHeight = 27

# --- SnippetID column (H) values were regenerated for this export. ---
$snippetIds = [ordered]@{
    "H2" = "mqrG"; "H3" = "jnBw"; "H4" = "qidl"; "H5" = "SUNv"; "H6" = "D84D";
    "H7" = "EoTA"; "H8" = "Gcb0"; "H9" = "Zeqp"; "H10" = "aWNi"; "H11" = "aWNi";
    "H12" = "aWNi"; "H13" = "aWNi"; "H14" = "Ivci"; "H15" = "ZZCv"; "H16" = "gQ3M";
    "H17" = "JpmJ"; "H18" = "JpmJ"; "H19" = "Bun9"; "H20" = "Fken"; "H21" = "QpWE";
    "H22" = "MQGj"; "H23" = "5hEN"; "H24" = "yV9X"; "H25" = "lfUb"; "H26" = "pZwi";
    "H27" = "qxob"; "H28" = "tc8C"; "H29" = "68Sc"; "H30" = "E5ZZ"; "H31" = "anrS";
    "H32" = "pIFe"; "H33" = "Vvby"; "H34" = "sQuy"; "H35" = "sQuy"; "H36" = "yZTD";
    "H37" = "yZTD"; "H38" = "tQ7q"; "H39" = "tQ7q"; "H40" = "Mtbr"; "H41" = "Mtbr";
    "H42" = "lMbn"; "H43" = "lMbn"; "H44" = "wRVG"; "H45" = "wRVG"; "H46" = "NVPq";
    "H47" = "NVPq"; "H48" = "LUnl"; "H49" = "hiPf"; "H50" = "hiPf"; "H51" = "ShN2";
    "H52" = "ShN2"; "H53" = "6hRb"; "H54" = "6hRb"; "H55" = "6hRb"; "H56" = "s5iM";
    "H57" = "s5iM"; "H58" = "cpha"; "H59" = "qeo2"; "H60" = "G4JV"; "H61" = "TFED"
}

foreach ($addr in $snippetIds.Keys) {
    $ws.Range($addr).Value = $snippetIds[$addr]
}
